$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PTON")

# Row 4 (Inventory) updates
$ws.Range("B4").Value = 523000000.0
$ws.Range("C4").Value = 364000000.0
$ws.Range("D4").Value = 245000000.0
$ws.Range("E4").Value = 194000000.0
$ws.Range("F4").Value = 244000000.0

# Row 13 (Accounts Payable) updates
$ws.Range("B13").Value = 721000000.0
$ws.Range("C13").Value = 275000000.0
$ws.Range("D13").Value = 136000000.0
$ws.Range("E13").Value = 142000000.0
$ws.Range("F13").Value = 183000000.0
